$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.251.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.569.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.557.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.207"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.606"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000288"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "694.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.131.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.302.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.549.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.69%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.926"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "587.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.667.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.144"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0761"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0438"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.345"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("E48").Value = "  +5.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
